$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the week's dates (row 5): shift the week forward by one day
$ws.Range("B5").Value = 44452
$ws.Range("C5").Value = 44453
$ws.Range("D5").Value = 44454
$ws.Range("E5").Value = 44455
$ws.Range("F5").Value = 44456
$ws.Range("G5").Value = 44457
$ws.Range("H5").Value = 44458

# Lecture row: moved from Monday (C6) to Sunday/Monday (B6)
$ws.Range("C6").ClearContents()
$ws.Range("B6").Value = 1

# Team Meeting row: moved from Wednesday (E8) to Tuesday (D8)
$ws.Range("E8").ClearContents()
$ws.Range("D8").Value = 1

# Sponsor Meeting row: moved from Monday (C9) to Sunday/Monday (B9)
$ws.Range("C9").ClearContents()
$ws.Range("B9").Value = 1

# Update selection
$ws.Range("K11").Select()
